# Fruta / hortaliza, semanal
# Insert one new weekly record as row 666 (pushing the existing rows 666-720
# down to 667-721), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 666:720 down one position, creating a blank (but format-carrying) row 666.
$ws.Rows("666:666").Insert()

# Populate the new row 666 with the new record's data.
$ws.Range("A666").Value = 9
$ws.Range("B666").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C666").Value = "Metropolitana"
$ws.Range("D666").Value = 44826
$ws.Range("E666").Value = 13
$ws.Range("F666").Value = 100112008
$ws.Range("G666").Value = "Coliflor"
$ws.Range("H666").Value = "Sin especificar"
$ws.Range("I666").Value = "Primera"
$ws.Range("J666").Value = 2200
$ws.Range("K666").Value = 800
$ws.Range("L666").Value = 1000
$ws.Range("M666").Value = 909
$ws.Range("N666").Value = '$/unidad'
$ws.Range("O666").Value = "Provincia de Melipilla"
$ws.Range("P666").Value = 909
$ws.Range("Q666").Value = 1
$ws.Range("R666").Value = "Hortaliza"
